# This script re-applies the weekly refresh of the "Alcachofa" price
# sheet: the data rows 2-26 (columns D, H, I, J, K, L, M, N, O, P, Q)
# get reshuffled, pulling each row's data from a different original row
# according to the mapping below (new row -> source row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of destination row number -> source row number (1-indexed sheet rows)
$rowMap = @{
    2  = 25
    3  = 2
    4  = 7
    5  = 5
    6  = 6
    7  = 11
    8  = 16
    9  = 26
    10 = 14
    11 = 15
    12 = 12
    13 = 21
    14 = 9
    15 = 13
    16 = 8
    17 = 10
    18 = 17
    19 = 18
    20 = 19
    21 = 20
    22 = 22
    23 = 4
    24 = 23
    25 = 3
    26 = 24
}

# Columns that participate in the reshuffle.
$cols = @("D", "H", "I", "J", "K", "L", "M", "N", "O", "P", "Q")

# Snapshot the current ("before") values for every row/column involved,
# so that reads are unaffected by writes we perform later.
$snapshot = @{}
foreach ($row in $rowMap.Values | Sort-Object -Unique) {
    $rowData = @{}
    foreach ($col in $cols) {
        $rowData[$col] = $ws.Range("$col$row").Value2
    }
    $snapshot[$row] = $rowData
}

# Write the snapshotted values into their destination rows.
foreach ($destRow in ($rowMap.Keys | Sort-Object)) {
    $srcRow = $rowMap[$destRow]
    $rowData = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value2 = $rowData[$col]
    }
}
